# 600 Essential Words.xlsx - add new vocabulary rows to the "Warranties" and
# "Bussines Planning" sheets, and move the active tab to "Bussines Planning".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Warranties sheet (3rd tab): 12 new words in A2:A13, plus a single blank-ish
# space value tucked away in H14.
# ---------------------------------------------------------------------------
$wsWarranties = $wb.Worksheets.Item("Warranties")

$warrantyWords = @(
    "characteristic",
    "consequence",
    "consider",
    "cover",
    "expiration",
    "frequently",
    "imply",
    "promise",
    "protect",
    "requtation",
    "require",
    "variety"
)

for ($i = 0; $i -lt $warrantyWords.Length; $i++) {
    $row = $i + 2
    $wsWarranties.Cells.Item($row, 1).Value = $warrantyWords[$i]
}

$wsWarranties.Range("H14").Value = " "

# Widen column A a bit to fit the new words.
$wsWarranties.Columns.Item(1).ColumnWidth = 14.71

# ---------------------------------------------------------------------------
# Bussines Planning sheet (4th tab): 12 new words in A2:A13.
# ---------------------------------------------------------------------------
$wsPlanning = $wb.Worksheets.Item("Bussines Planning")

$planningWords = @(
    "address",
    "avoid",
    "demonstrate",
    "develop",
    "evaluate",
    "gather",
    "offer",
    "primarily",
    "risk",
    "strategy",
    "strong",
    "substitution"
)

for ($i = 0; $i -lt $planningWords.Length; $i++) {
    $row = $i + 2
    $wsPlanning.Cells.Item($row, 1).Value = $planningWords[$i]
}

# ---------------------------------------------------------------------------
# Update selections to match where the author last left off, then make
# "Bussines Planning" the active tab.
# ---------------------------------------------------------------------------
[void]$wsWarranties.Range("H14").Select()
[void]$wsPlanning.Range("D14").Select()
[void]$wsPlanning.Activate()
